$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5 text changes from "Period from Freq." to "Desired Period"
$ws.Range("A5").Value = "Desired Period"

# Update selection to I11
$ws.Range("I11").Select()
